# Update the "Data Kategori" template worksheet:
#  - Replace the 4 header columns (No / Id Kategori / Kode Kategori / Nama Kategori)
#    with 2 headers (Kategori_Kode / Kategori_Nama)
#  - Clear the now-unused columns C:D
#  - Add an empty (unbolded) data row below the header
#  - Resize column A to fit the new header text
#  - Move the active selection down to B10:B11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous 4-column header row and any leftover rows/cols entirely
$ws.Cells.Clear()

# New 2-column header row (bold, same style as before: fontId 1 / s="1")
$ws.Range("A1").Value = "Kategori_Kode"
$ws.Range("B1").Value = "Kategori_Nama"
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").Font.Name = "Calibri"
$ws.Range("A1:B1").Font.Size = 11

# New empty (non-bold) data row under the header
$ws.Range("A2:B2").Font.Bold = $false
$ws.Range("A2:B2").Font.Name = "Calibri"
$ws.Range("A2:B2").Font.Size = 11

# Column A needs to be wider to fit the new header text (columns B:D keep
# their existing widths from the original template)
$ws.Columns.Item(1).ColumnWidth = 16.140625

# Move the selection like the saved file shows
$ws.Range("B10:B11").Select()
